$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells keep their text formatting (values like "0.4813",
# "29.448.15" must not be auto-converted to numbers/dates by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.448.15'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '1.908.85'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.68%  '
$ws.Range('D5').Value = '325.23'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').Value = '0.4813'
$ws.Range('E7').Value = '  +1.73%  '
$ws.Range('D8').Value = '0.4066'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.08160'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('D11').Value = '23.45'
$ws.Range('E11').Value = '  +3.53%  '
$ws.Range('D12').Value = '1.940.88'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').Value = '6.011'
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').Value = '7.168'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').Value = '90.31'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = '0.06794'
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '0.00001036'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').Value = '17.69'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').Value = '29.466.04'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('D22').Value = '5.616'
$ws.Range('E22').Value = '  +1.76%  '
$ws.Range('D23').Value = '11.71'
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('D25').Value = '2.136.65'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').Value = '156.27'
$ws.Range('E26').Value = '  +1.58%  '
$ws.Range('D27').Value = '6.422'
$ws.Range('E27').Value = '  +6.34%  '
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '120.18'
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').Value = '1.025'
$ws.Range('E31').Value = '  -4.58%  '
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('D33').Value = '5.522'
$ws.Range('E33').Value = '  +2.52%  '
$ws.Range('D34').Value = '3.559'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('E35').Value = '  -2.40%  '
$ws.Range('D36').Value = '0.02269'
$ws.Range('E36').Value = '  +0.98%  '
$ws.Range('D37').Value = '0.06099'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').Value = '1.176'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('D39').Value = '10.85'
$ws.Range('E39').Value = '  +7.27%  '
$ws.Range('D40').Value = '0.5961'
$ws.Range('E40').Value = '  +1.78%  '
$ws.Range('D41').Value = '7.984'
$ws.Range('E41').Value = '  -3.19%  '
$ws.Range('D42').Value = '0.1855'
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('D43').Value = '1.280'
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').Value = '2.392'
$ws.Range('E44').Value = '  -4.95%  '
$ws.Range('D45').Value = '12.55'
$ws.Range('E45').Value = '  +3.42%  '
$ws.Range('D46').Value = '0.07595'
$ws.Range('E46').Value = '  -3.73%  '
$ws.Range('D47').Value = '0.5570'
$ws.Range('E47').Value = '  +1.05%  '
$ws.Range('D48').Value = '1.944'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').Value = '116.31'
$ws.Range('E49').Value = '  +2.65%  '
$ws.Range('D50').Value = '72.56'
$ws.Range('E50').Value = '  +1.68%  '
$ws.Range('D51').Value = '2.407'
$ws.Range('E51').Value = '  +2.53%  '
